$d = $word.ActiveDocument

$end = $d.Content
$end.Collapse(0)  # wdCollapseEnd
$end.InsertParagraphAfter()
$end.Collapse(0)
$end.Move(4, 1) | Out-Null  # wdParagraph, move into new paragraph

$p = $d.Paragraphs.Last
$p.Range.Text = "Implementado el diagrama 2D, pero la librería no está agregada en el JAR"
$p.Style = "Cita"
$p.Alignment = 2  # wdAlignParagraphRight
